$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2026.7778
$ws.Range("I15").Value = 2026.7778
$ws.Range("K15").Value = 6080.3334
$ws.Range("M15").Value = -5911.3334
$ws.Range("H19").Value = 663.5
$ws.Range("I19").Value = 755.8570999999999
$ws.Range("J19").Value = 571.1429000000001
$ws.Range("K19").Value = 755.8570999999999
$ws.Range("L19").Value = 571.1429000000001
$ws.Range("M19").Value = -580.8570999999999
$ws.Range("N19").Value = -921.1429000000001
$ws.Range("H58").Value = 1789.6666
$ws.Range("J58").Value = 4999.5
$ws.Range("L58").Value = 14998.5
$ws.Range("N58").Value = -15298.5
$ws.Range("H70").Value = 1933.826
$ws.Range("J70").Value = 3180.7778
$ws.Range("L70").Value = 9542.3334
$ws.Range("N70").Value = -10082.3334
$ws.Range("H73").Value = 1933.826
$ws.Range("J73").Value = 3180.7778
$ws.Range("L73").Value = 9542.3334
$ws.Range("N73").Value = -11414.3334
$ws.Range("H88").Value = 6375
$ws.Range("I88").Value = 20000
$ws.Range("J88").Value = 1833.3334
$ws.Range("K88").Value = 20000
$ws.Range("L88").Value = 1833.3334
$ws.Range("M88").Value = -19594
$ws.Range("N88").Value = -2645.3334
$ws.Range("H91").Value = 6375
$ws.Range("I91").Value = 20000
$ws.Range("J91").Value = 1833.3334
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 1833.3334
$ws.Range("M91").Value = -18596
$ws.Range("N91").Value = -4641.3334
$ws.Range("H93").Value = 34500
$ws.Range("J93").Value = 34500
$ws.Range("L93").Value = 34500
$ws.Range("N93").Value = -39492
$ws.Range("H132").Value = 776.37256
$ws.Range("I132").Value = 789.36365
$ws.Range("K132").Value = 2368.09095
$ws.Range("M132").Value = 161.9090500000002
$ws.Range("H137").Value = 2613.5217
$ws.Range("I137").Value = 1917.0769
$ws.Range("J137").Value = 3518.9
$ws.Range("K137").Value = 5751.2307
$ws.Range("L137").Value = 10556.7
$ws.Range("M137").Value = -3201.2307
$ws.Range("N137").Value = -15656.7
$ws.Range("H141").Value = 2447.625
$ws.Range("I141").Value = 2495.282
$ws.Range("K141").Value = 7485.846
$ws.Range("M141").Value = -2305.846

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 21536.334
$ws.Range("I31").Value = 17690
$ws.Range("K31").Value = 17690
$ws.Range("M31").Value = -17396
$ws.Range("H56").Value = 339666.66
$ws.Range("I56").Value = 9500
$ws.Range("J56").Value = 1000000
$ws.Range("K56").Value = 9500
$ws.Range("L56").Value = 1000000
$ws.Range("M56").Value = -8758
$ws.Range("N56").Value = -1001484
$ws.Range("H110").Value = 3616.5
$ws.Range("I110").Value = 3533.15
$ws.Range("J110").Value = 4450
$ws.Range("K110").Value = 3533.15
$ws.Range("L110").Value = 4450
$ws.Range("M110").Value = -1488.15
$ws.Range("N110").Value = -8540
$ws.Range("H122").Value = 3446.3928
$ws.Range("I122").Value = 2504.762
$ws.Range("K122").Value = 7514.286
$ws.Range("M122").Value = -5064.286

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2918.3333
$ws.Range("I86").Value = 2918.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2918.3333
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1795.3333
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2918.3333
$ws.Range("I89").Value = 2918.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 14591.6665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -8975.666499999999
$ws.Range("N89").ClearContents()
$ws.Range("H102").Value = 26872.584
$ws.Range("I102").Value = 25679.182
$ws.Range("K102").Value = 25679.182
$ws.Range("M102").Value = -22434.182

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4444.5
$ws.Range("J58").Value = 5333.3335
$ws.Range("L58").Value = 5333.3335
$ws.Range("N58").Value = -5739.3335
$ws.Range("H132").Value = 6152.7
$ws.Range("I132").Value = 6533.3335
$ws.Range("K132").Value = 19600.0005
$ws.Range("M132").Value = -17070.0005
$ws.Range("H136").Value = 4444.5
$ws.Range("J136").Value = 5333.3335
$ws.Range("L136").Value = 16000.0005
$ws.Range("N136").Value = -21100.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2053
$ws.Range("I129").Value = 1320
$ws.Range("J129").Value = 2252.9092
$ws.Range("K129").Value = 3960
$ws.Range("L129").Value = 6758.7276
$ws.Range("M129").Value = 1040
$ws.Range("N129").Value = -16758.7276

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6251.7114
$ws.Range("I46").Value = 1179.45
$ws.Range("K46").Value = 1179.45
$ws.Range("M46").Value = -991.45
$ws.Range("H61").Value = 3269.8076
$ws.Range("I61").Value = 2195.1765
$ws.Range("K61").Value = 2195.1765
$ws.Range("M61").Value = -1993.1765
$ws.Range("H97").Value = 15326
$ws.Range("J97").Value = 15326
$ws.Range("L97").Value = 15326
$ws.Range("N97").Value = -17308
$ws.Range("H113").Value = 3269.8076
$ws.Range("I113").Value = 2195.1765
$ws.Range("K113").Value = 2195.1765
$ws.Range("M113").Value = -25.17650000000003
$ws.Range("H122").Value = 8259.666999999999
$ws.Range("I122").Value = 7201.8335
$ws.Range("J122").Value = 9317.5
$ws.Range("K122").Value = 21605.5005
$ws.Range("L122").Value = 27952.5
$ws.Range("M122").Value = -19155.5005
$ws.Range("N122").Value = -32852.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 44624.5
$ws.Range("J63").Value = 44624.5
$ws.Range("L63").Value = 44624.5
$ws.Range("N63").Value = -45872.5
$ws.Range("H66").Value = 44624.5
$ws.Range("J66").Value = 44624.5
$ws.Range("L66").Value = 133873.5
$ws.Range("N66").Value = -140113.5
$ws.Range("H101").Value = 19499.75
$ws.Range("J101").Value = 19499.75
$ws.Range("L101").Value = 19499.75
$ws.Range("N101").Value = -25989.75
$ws.Range("H104").Value = 26500
$ws.Range("J104").Value = 26500
$ws.Range("L104").Value = 26500
$ws.Range("N104").Value = -33488
$ws.Range("H132").Value = 2724.6924
$ws.Range("I132").Value = 1618.5555
$ws.Range("K132").Value = 4855.666499999999
$ws.Range("M132").Value = -2325.666499999999
$ws.Range("H136").Value = 12503712
$ws.Range("I136").Value = 20003170
$ws.Range("J136").Value = 4616.933
$ws.Range("K136").Value = 60009510
$ws.Range("L136").Value = 13850.799
$ws.Range("M136").Value = -60006960
$ws.Range("N136").Value = -18950.799

Write-Output "Applied all profit sheet updates"